# Apply the "basis.xlsx" data-file update:
#   - rename German headers to English headers
#   - drop the bold styling from the "Customer Number" / ID column header and
#     give the whole column a consistent (non-bold) font
#   - give the other three columns (First Name / Age / City) an explicit,
#     consistent font as well
#   - autofit column A to the new, longer header text
#   - move the active cell/selection to E11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text: German -> English -------------------------------------
$ws.Range("A1").Value = "Customer Number"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "City"

# --- 2. Fonts / styling ------------------------------------------------------
# The header cell (A1) loses the bold weight it used to have, so the whole
# id column (header included) now shares one consistent look; then the
# font name gets normalized for the whole column.
$ws.Range("A1").Font.Bold = $false
$colA = $ws.Range("A1:A11")
$colA.Font.Name = "Calibri"

# Columns B:D (including their header row) get an explicit font matching the
# workbook's default look.
$colBD = $ws.Range("B1:D11")
$colBD.Font.Name = "Calibri"

# --- 3. Column width ----------------------------------------------------------
# Autofit column A for the new, wider header ("Customer Number"), then settle
# on the final width used by the workbook.
$col = $ws.Columns.Item(1)
$col.AutoFit()
$col.ColumnWidth = 16.6

# --- 4. Selection -------------------------------------------------------------
$ws.Range("E11").Select()
